$d = $word.ActiveDocument

# Locate the trailing " Absolute Deviation Accuracy: 84.3% Total Deviation
# Accuracy: 99.4%" text that must be removed from the end of the run, and
# replaced by two new runs: one holding "." and one holding a single
# trailing space - all three runs sharing identical run formatting.
$target = $d.Content
$found = $target.Find.Execute(" Absolute Deviation Accuracy: 84.3% Total Deviation Accuracy: 99.4%")

if ($found) {
    $insertPos = $target.Start

    # Remove the unwanted trailing text from the original run.
    $target.Text = ""

    # Insert the replacement ". " as a simple continuation of the run that
    # precedes it, so it fully inherits that run's character formatting
    # (font, color, size, etc.) instead of falling back to any default.
    $tail = $d.Range($insertPos, $insertPos)
    $tail.InsertAfter(". ")

    # Force the "." to live in its own <w:r> (identical formatting to its
    # neighbours) by toggling a property on just that character and back
    # to its original value - this breaks run-coalescing without changing
    # the resulting visible/resolved formatting at all.
    $dotRange = $d.Range($insertPos, $insertPos + 1)
    $dotRange.Bold = 1
    $dotRange.Bold = 0
}
